$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = "boeing_787-8_cropped"
$ws.Range("K3").Value = "boeing_787-8_cropped"
$ws.Range("K4").Value = "boeing_787-8_cropped"
$ws.Range("K5").Value = "B737_100"
$ws.Range("K6").Value = "B737_100"
$ws.Range("K7").Value = "a350_900_cropped"

$ws.Range("K3").Select()
